# Normalize the "Recorded By" column (G): when the literal entry "System"
# appears anywhere in the comma-separated list of recorders, move it to the
# front of the list (list order is otherwise preserved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Move-SystemFirst($val) {
    if ($val -eq $null) { return $val }
    if ($val -notlike "*System*") { return $val }

    $parts = $val -split ",\s*"
    $idx = [array]::IndexOf($parts, "System")
    if ($idx -le 0) {
        # "System" absent, or already first -> nothing to do.
        return $val
    }

    $newParts = New-Object System.Collections.ArrayList
    [void]$newParts.Add("System")
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $idx) {
            [void]$newParts.Add($parts[$i])
        }
    }
    return ($newParts -join ", ")
}

$dims = $ws.UsedRange
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    $newVal = Move-SystemFirst $val
    # NOTE: use .Equals() (ordinal / case-sensitive) rather than -eq/-ne,
    # since -eq on strings here is case-insensitive and would wrongly treat
    # "system, System, x" and "System, system, x" as identical.
    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
